$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 109.833336
$ws.Range("I38").Value = 109.833336
$ws.Range("K38").Value = 329.500008
$ws.Range("M38").Value = 42.49999200000002
$ws.Range("H98").Value = 722.55554
$ws.Range("I98").Value = 847.0625
$ws.Range("K98").Value = 847.0625
$ws.Range("M98").Value = 650.9375
$ws.Range("H100").Value = 3772.5173
$ws.Range("I100").Value = 1186.8462
$ws.Range("J100").Value = 5873.375
$ws.Range("K100").Value = 1186.8462
$ws.Range("L100").Value = 5873.375
$ws.Range("M100").Value = -645.8462
$ws.Range("N100").Value = -6955.375
$ws.Range("H122").Value = 722.55554
$ws.Range("I122").Value = 847.0625
$ws.Range("K122").Value = 2541.1875
$ws.Range("M122").Value = -91.1875
$ws.Range("H132").Value = 14335.6
$ws.Range("I132").Value = 15302.714
$ws.Range("K132").Value = 45908.142
$ws.Range("M132").Value = -43378.142

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 500225
$ws.Range("I16").Value = 500225
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 500225
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -499938
$ws.Range("N16").ClearContents()
$ws.Range("H19").Value = 1008
$ws.Range("I19").Value = 1008
$ws.Range("K19").Value = 1008
$ws.Range("M19").Value = -779
$ws.Range("H32").Value = 197652.08
$ws.Range("I32").Value = 203478.16
$ws.Range("K32").Value = 203478.16
$ws.Range("M32").Value = -203191.16
$ws.Range("H45").Value = 2684.3
$ws.Range("I45").Value = 1834.7142
$ws.Range("J45").Value = 4666.6665
$ws.Range("K45").Value = 1834.7142
$ws.Range("L45").Value = 4666.6665
$ws.Range("M45").Value = -1457.7142
$ws.Range("N45").Value = -5420.6665
$ws.Range("H88").Value = 2561.0386
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 2623.48
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 2623.48
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -3435.48
$ws.Range("H91").Value = 2561.0386
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 2623.48
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 2623.48
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -5431.48
$ws.Range("H122").Value = 2675.8235
$ws.Range("I122").Value = 2159.16
$ws.Range("J122").Value = 4111
$ws.Range("K122").Value = 6477.48
$ws.Range("L122").Value = 12333
$ws.Range("M122").Value = -4027.48
$ws.Range("N122").Value = -17233
$ws.Range("H133").Value = 67450
$ws.Range("J133").Value = 67450
$ws.Range("L133").Value = 67450
$ws.Range("N133").Value = -72510
$ws.Range("H134").Value = 104798.4
$ws.Range("J134").Value = 104798.4
$ws.Range("L134").Value = 104798.4
$ws.Range("N134").Value = -114938.4

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6460.5713
$ws.Range("I134").Value = 6704.3887
$ws.Range("K134").Value = 20113.1661
$ws.Range("M134").Value = -17578.1661
$ws.Range("H135").Value = 66491
$ws.Range("J135").Value = 66491
$ws.Range("L135").Value = 66491
$ws.Range("N135").Value = -76631

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 129863.25
$ws.Range("I16").Value = 7181.2
$ws.Range("K16").Value = 7181.2
$ws.Range("M16").Value = -6894.2
$ws.Range("H19").Value = 775.75
$ws.Range("I19").Value = 775.75
$ws.Range("K19").Value = 775.75
$ws.Range("M19").Value = -605.75
$ws.Range("H24").Value = 775.75
$ws.Range("I24").Value = 775.75
$ws.Range("K24").Value = 775.75
$ws.Range("M24").Value = -605.75
$ws.Range("H31").Value = 2124.9333
$ws.Range("I31").Value = 2124.9333
$ws.Range("K31").Value = 2124.9333
$ws.Range("M31").Value = -1829.9333
$ws.Range("H34").Value = 2124.9333
$ws.Range("I34").Value = 2124.9333
$ws.Range("K34").Value = 2124.9333
$ws.Range("M34").Value = -1922.9333
$ws.Range("H62").Value = 4057
$ws.Range("I62").Value = 3680
$ws.Range("K62").Value = 3680
$ws.Range("M62").Value = -3056
$ws.Range("H65").Value = 4057
$ws.Range("I65").Value = 3680
$ws.Range("K65").Value = 18400
$ws.Range("M65").Value = -15280
$ws.Range("H99").Value = 27426
$ws.Range("J99").Value = 3500
$ws.Range("L99").Value = 3500
$ws.Range("N99").Value = -6496
$ws.Range("H113").Value = 129863.25
$ws.Range("I113").Value = 7181.2
$ws.Range("K113").Value = 7181.2
$ws.Range("M113").Value = -5011.2
$ws.Range("H126").Value = 27426
$ws.Range("J126").Value = 3500
$ws.Range("L126").Value = 10500
$ws.Range("N126").Value = -15440

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 799
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H12").Value = 1093
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1093
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 3279
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -3625
$ws.Range("H131").Value = 12000
$ws.Range("J131").Value = 12000
$ws.Range("L131").Value = 36000
$ws.Range("N131").Value = -46080

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 12000
$ws.Range("J18").Value = 12000
$ws.Range("L18").Value = 12000
$ws.Range("N18").Value = -12586
$ws.Range("H21").Value = 4599.3335
$ws.Range("I21").Value = 4599.3335
$ws.Range("K21").Value = 4599.3335
$ws.Range("M21").Value = -4426.3335
$ws.Range("H30").Value = 4599.3335
$ws.Range("I30").Value = 4599.3335
$ws.Range("K30").Value = 4599.3335
$ws.Range("M30").Value = -4494.3335
$ws.Range("H102").Value = 2842.2144
$ws.Range("I102").Value = 2945.4614
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 2945.4614
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -1323.4614
$ws.Range("N102").Value = -4744
$ws.Range("H122").Value = 2327.7144
$ws.Range("I122").Value = 1526.8572
$ws.Range("J122").Value = 3128.5715
$ws.Range("K122").Value = 4580.571599999999
$ws.Range("L122").Value = 9385.7145
$ws.Range("M122").Value = -2130.571599999999
$ws.Range("N122").Value = -14285.7145
$ws.Range("H132").Value = 12269.774
$ws.Range("I132").Value = 12615.434
$ws.Range("K132").Value = 37846.302
$ws.Range("M132").Value = -35316.302

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2533.9
$ws.Range("I22").Value = 1945.6923
$ws.Range("J22").Value = 3626.2856
$ws.Range("K22").Value = 1945.6923
$ws.Range("L22").Value = 3626.2856
$ws.Range("M22").Value = -1650.6923
$ws.Range("N22").Value = -4216.2856
$ws.Range("H27").Value = 2533.9
$ws.Range("I27").Value = 1945.6923
$ws.Range("J27").Value = 3626.2856
$ws.Range("K27").Value = 1945.6923
$ws.Range("L27").Value = 3626.2856
$ws.Range("M27").Value = -1838.6923
$ws.Range("N27").Value = -3840.2856
$ws.Range("H46").Value = 4077.95
$ws.Range("I46").Value = 1113.6
$ws.Range("J46").Value = 5066.067
$ws.Range("K46").Value = 1113.6
$ws.Range("L46").Value = 5066.067
$ws.Range("M46").Value = -925.5999999999999
$ws.Range("N46").Value = -5442.067
$ws.Range("H132").Value = 2730.3103
$ws.Range("J132").Value = 3897.2307
$ws.Range("L132").Value = 11691.6921
$ws.Range("N132").Value = -16751.6921
$ws.Range("H133").Value = 88777
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2697.6365
$ws.Range("I132").Value = 2052.8667
$ws.Range("K132").Value = 6158.6001
$ws.Range("M132").Value = -3628.6001
